$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 258 - this shifts the former rows 258..306
# down to become rows 259..307, preserving their content untouched.
$ws.Rows.Item(258).Insert()

# Populate the newly inserted row 258 with the new record
# (Feria Lagunitas de Puerto Montt, Durazno, Early Majestic).
$ws.Cells.Item(258, 1).Value = 4
$ws.Cells.Item(258, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(258, 3).Value = "Los Lagos"
$ws.Cells.Item(258, 4).Value = 45258
$ws.Cells.Item(258, 5).Value = 10
$ws.Cells.Item(258, 6).Value = "Fruta"
$ws.Cells.Item(258, 7).Value = 100103
$ws.Cells.Item(258, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(258, 9).Value = 100103004
$ws.Cells.Item(258, 10).Value = "Durazno"
$ws.Cells.Item(258, 11).Value = "Early Majestic"
$ws.Cells.Item(258, 12).Value = "Especial"
$ws.Cells.Item(258, 13).Value = 300
$ws.Cells.Item(258, 14).Value = 30000
$ws.Cells.Item(258, 15).Value = 30000
$ws.Cells.Item(258, 16).Value = 30000
$ws.Cells.Item(258, 17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(258, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(258, 19).Value = 2143
$ws.Cells.Item(258, 20).Value = 14
